$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.543.01"
$ws.Range("D3").Value = "3.503.47"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.80"
$ws.Range("E5").Value = "  -1.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "194.76"
$ws.Range("E6").Value = "  +3.04%  "
$ws.Range("E7").Value = "  -0.46%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  -4.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.644"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.11"
$ws.Range("E11").Value = "  +0.70%  "
$ws.Range("E12").Value = "  -2.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.45"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "4.055.37"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "594.50"
$ws.Range("E15").Value = "  -2.73%  "
$ws.Range("D16").Value = "69.685.87"
$ws.Range("E16").Value = "  +0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.00"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.65"
$ws.Range("E18").Value = "  +1.43%  "
$ws.Range("E19").Value = "  +2.44%  "
$ws.Range("D20").Value = "3.496.17"
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "18.11"
$ws.Range("E22").Value = "  +6.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.26"
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.72"
$ws.Range("E24").Value = "  -4.48%  "
$ws.Range("E25").Value = "  -1.37%  "
$ws.Range("E26").Value = "  +4.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  -1.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.48"
$ws.Range("E28").Value = "  -1.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.04"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.29"
$ws.Range("E30").Value = "  +9.04%  "
$ws.Range("E31").Value = "  +1.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.33"
$ws.Range("E32").Value = "  -1.43%  "
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.11"
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "3.728.98"
$ws.Range("E35").Value = "  +3.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.12"
$ws.Range("E36").Value = "  +0.37%  "
$ws.Range("D37").Value = "0.0₃0807"
$ws.Range("E37").Value = "  +5.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("E40").Value = "  -0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.20"
$ws.Range("E41").Value = "  -0.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "486.70"
$ws.Range("E42").Value = "  -3.47%  "
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  -3.00%  "
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.38"
$ws.Range("E49").Value = "  -3.60%  "
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("E51").Value = "  +10.13%  "
